$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Maps 2 RGB")

# --- Log metrics for rows 186-193 (versions V9.07.5 .. V9.07.12 worth of PSNR/SSIM columns B:I) ---
$ws.Range("B186").Value = 15.4708
$ws.Range("C186").Value = 0.8209
$ws.Range("D186").Value = 8.3193
$ws.Range("E186").Value = 0.7242
$ws.Range("F186").Value = 15.119
$ws.Range("G186").Value = 0.8018
$ws.Range("H186").Value = 21.3439
$ws.Range("I186").Value = 0.7992

$ws.Range("B187").Value = 13.1619
$ws.Range("C187").Value = 0.8035
$ws.Range("D187").Value = 8.5848
$ws.Range("E187").Value = 0.7245
$ws.Range("F187").Value = 16.0564
$ws.Range("G187").Value = 0.8154
$ws.Range("H187").Value = 21.7594
$ws.Range("I187").Value = 0.7965

$ws.Range("B188").Value = 13.1266
$ws.Range("C188").Value = 0.8003
$ws.Range("D188").Value = 8.7202
$ws.Range("E188").Value = 0.7289
$ws.Range("F188").Value = 16.1456
$ws.Range("G188").Value = 0.8202
$ws.Range("H188").Value = 17.9141
$ws.Range("I188").Value = 0.7889

$ws.Range("B189").Value = 13.1436
$ws.Range("C189").Value = 0.8015
$ws.Range("D189").Value = 8.1607
$ws.Range("E189").Value = 0.7306
$ws.Range("F189").Value = 16.2872
$ws.Range("G189").Value = 0.832
$ws.Range("H189").Value = 17.0889
$ws.Range("I189").Value = 0.7868

$ws.Range("B190").Value = 13.1805
$ws.Range("C190").Value = 0.8008
$ws.Range("D190").Value = 8.3765
$ws.Range("E190").Value = 0.7283
$ws.Range("F190").Value = 15.6667
$ws.Range("G190").Value = 0.8167
$ws.Range("H190").Value = 17.7053
$ws.Range("I190").Value = 0.7917

$ws.Range("B191").Value = 13.0721
$ws.Range("C191").Value = 0.7992
$ws.Range("D191").Value = 8.121
$ws.Range("E191").Value = 0.7328
$ws.Range("F191").Value = 16.3556
$ws.Range("G191").Value = 0.8303
$ws.Range("H191").Value = 16.6251
$ws.Range("I191").Value = 0.7817

$ws.Range("B192").Value = 13.2451
$ws.Range("C192").Value = 0.8008
$ws.Range("D192").Value = 8.2849
$ws.Range("E192").Value = 0.7303
$ws.Range("F192").Value = 14.8488
$ws.Range("G192").Value = 0.7966
$ws.Range("H192").Value = 18.8231
$ws.Range("I192").Value = 0.7935

$ws.Range("B193").Value = 13.4854
$ws.Range("C193").Value = 0.8007
$ws.Range("D193").Value = 8.7506
$ws.Range("E193").Value = 0.7306
$ws.Range("F193").Value = 16.5245
$ws.Range("G193").Value = 0.818
$ws.Range("H193").Value = 18.2232
$ws.Range("I193").Value = 0.7798

# --- Add new version labels for rows 198-211 (column A) ---
$ws.Range("A198").Value = "V9.07.5"
$ws.Range("A199").Value = "V9.07.6"
$ws.Range("A200").Value = "V9.07.7"
$ws.Range("A201").Value = "V9.07.8"
$ws.Range("A202").Value = "V9.07.9"
$ws.Range("A203").Value = "V9.07.10"
$ws.Range("A204").Value = "V9.07.11"
$ws.Range("A205").Value = "V9.07.12"
$ws.Range("A206").Value = "V9.07.13"
$ws.Range("A207").Value = "V9.07.14"
$ws.Range("A208").Value = "V9.07.15"
$ws.Range("A209").Value = "V9.07.16"
$ws.Range("A210").Value = "V9.07.17"
$ws.Range("A211").Value = "V9.07.18"

# --- Reflect the final selection state on the sheet ---
$ws.Activate()
$ws.Range("I188").Select()
